# Implement csv module error handling
# Appends one new data row to each of the four log sheets, mirroring the
# existing row layout/format (column A keeps the datetime number format of
# the row above; columns B-E are textual hex byte dumps; F-I are numeric).

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $row, $timeValue, $totalHex, $idHex, $actualHex, $checksumHex, $totalDec, $idDec, $actualDec, $checksumDec)

    $prevRow = $row - 1

    $ws.Cells.Item($row, 1).Value = $timeValue
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    $ws.Cells.Item($row, 2).Value = $totalHex
    $ws.Cells.Item($row, 3).Value = $idHex
    $ws.Cells.Item($row, 4).Value = $actualHex
    $ws.Cells.Item($row, 5).Value = $checksumHex

    $ws.Cells.Item($row, 6).Value = $totalDec

    if ($idDec -is [string]) {
        # Large ID_DEC values lose precision as a double in this sheet, so
        # the source data keeps them as plain text - match that.
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $idDec
    } else {
        $ws.Cells.Item($row, 7).Value = $idDec
    }

    $ws.Cells.Item($row, 8).Value = $actualDec
    $ws.Cells.Item($row, 9).Value = $checksumDec
}

# --- Sheet "ROW50-FE-LIFTER": new row 52 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-LogRow $ws1 52 45750.19417305556 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x62" "0xe" 400 568631262647114000000000.0 354 14

# --- Sheet "ROW50-MID-LIFTER": new row 54 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-LogRow $ws2 54 45750.16596064815 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x66" "0x19" 400 "568631262647113771663628" 358 25

# --- Sheet "ROW11-FE-LIFTER": new row 52 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-LogRow $ws3 52 45750.22506553241 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x62" "0x14" 400 568631262647114000000000.0 354 20

# --- Sheet "ROW11-MID-LIFTER": new row 52 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-LogRow $ws4 52 45750.35911744213 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x6a" "0x19" 400 568631262647114000000000.0 362 25
